$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Zipcode")

# Fix the incorrect zipcode value for Columbus, OH (B4): 432215 -> 43215
$ws.Range("B4").Value = 43215

# Update the active selection to B4 to match the saved cursor position
$ws.Range("B4").Select()
